$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S21").Value = 38646

$ws.Range("S22").Value = 38253
$ws.Range("V22").ClearContents()
$ws.Range("W22").ClearContents()

$ws.Range("S23").Value = 38679

$ws.Range("S24").Value = 38098
$ws.Range("V24").Value = 3166.905916666667
$ws.Range("W24").Value = 8.356324842185151
